$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.254.91"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.852.83"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  -0.53%  "
$ws.Range("D5").Value = "'314.14"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "'0.4606"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("D8").Value = "'0.3706"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "'0.07299"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").Value = "'0.8873"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").Value = "'20.09"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").Value = "'0.07806"
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("D13").Value = "1.891.46"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "'5.388"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "'6.528"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "'91.42"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'1.002"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "'0.000008938"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "'14.77"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "27.277.85"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("D22").Value = "'5.110"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "2.048.35"
$ws.Range("E24").Value = "  -5.77%  "
$ws.Range("E25").Value = "  +4.98%  "
$ws.Range("D26").Value = "'151.84"
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").Value = "'18.45"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").Value = "'2.058"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").Value = "'115.88"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").Value = "'5.063"
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("D31").Value = "'0.08827"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("D32").Value = "'3.098"
$ws.Range("E32").Value = "  +4.54%  "
$ws.Range("D33").Value = "'0.7707"
$ws.Range("E33").Value = "  +5.68%  "
$ws.Range("E34").Value = "  +3.38%  "
$ws.Range("D35").Value = "'4.502"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("D36").Value = "'2.754"
$ws.Range("E36").Value = "  +11.90%  "
$ws.Range("D37").Value = "'1.084"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").Value = "'0.01950"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "'0.05261"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").Value = "'2.947"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").Value = "'7.072"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").Value = "'0.5122"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Value = "'0.1634"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "'8.390"
$ws.Range("E44").Value = "  +2.56%  "
$ws.Range("D45").Value = "'0.4799"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").Value = "'10.31"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("D49").Value = "'1.642"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "'0.06213"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "'65.67"
$ws.Range("E51").Value = "  +0.84%  "
